$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '269.49'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '22.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '6.333'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06181'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.387'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8294'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1606'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08366'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03397'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.03193'
$ws.Range('D14').Style = 'Normal'
$ws.Range('B15').NumberFormat = '@'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').NumberFormat = '@'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.09324'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.839'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.001643'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('E17').Style = 'Normal'
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('B18').Style = 'Normal'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('C18').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.04725'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '17CoinExTokenCET'
$ws.Range('E18').Style = 'Normal'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'TigerCash'
$ws.Range('B19').Style = 'Normal'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('C19').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.006328'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '18TigerCashTCH'
$ws.Range('E19').Style = 'Normal'
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'HotbitToken'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.005672'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '19HotbitTokenHTB'
$ws.Range('E20').Style = 'Normal'
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'BitKan'
$ws.Range('B21').Style = 'Normal'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('C21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.001076'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '20BitKanKAN'
$ws.Range('E21').Style = 'Normal'
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'NitroEx'
$ws.Range('B22').Style = 'Normal'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('C22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0001499'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '21NitroExNTX'
$ws.Range('E22').Style = 'Normal'
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'LEO'
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('C23').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.725'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '22LEOLEO'
$ws.Range('E23').Style = 'Normal'
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'BTSEToken'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.324'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '23BTSETokenBTSE'
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'BitpandaEcosystemToken'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.3308'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '24BitpandaEcosystemTokenBEST'
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'ProBitToken'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1239'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '25ProBitTokenPROB'
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0002701'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04705'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.006948'
$ws.Range('D41').Style = 'Normal'
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'CEJI'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.003796'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '41CEJICEJIWorstin24h'
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'BKEXToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1159'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '42BKEXTokenBKK'
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01190'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00006262'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0009896'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000749'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9191'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002267'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00001399'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '49CryptobidCoinCBC'
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.01239'
$ws.Range('D51').Style = 'Normal'
